$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '64.161.95'
$c.Style = $origStyle
$ws.Range('E2').Value = '  +1.02%  '
$c = $ws.Range('D3')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.127.43'
$c.Style = $origStyle
$ws.Range('E3').Value = '  +0.75%  '
$ws.Range('E4').Value = '  +0.12%  '
$c = $ws.Range('D5')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '601.94'
$c.Style = $origStyle
$ws.Range('E5').Value = '  -0.58%  '
$c = $ws.Range('D6')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '142.69'
$c.Style = $origStyle
$ws.Range('E6').Value = '  -0.61%  '
$ws.Range('E7').Value = '  -0.20%  '
$c = $ws.Range('D8')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.122.62'
$c.Style = $origStyle
$ws.Range('E8').Value = '  +0.66%  '
$c = $ws.Range('D9')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.521'
$c.Style = $origStyle
$ws.Range('E9').Value = '  +1.11%  '
$c = $ws.Range('D10')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.149'
$c.Style = $origStyle
$ws.Range('E10').Value = '  +0.15%  '
$c = $ws.Range('D11')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '5.35'
$c.Style = $origStyle
$ws.Range('E11').Value = '  +2.19%  '
$c = $ws.Range('D12')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.468'
$c.Style = $origStyle
$ws.Range('E12').Value = '  +0.53%  '
$c = $ws.Range('D13')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.0000254'
$c.Style = $origStyle
$ws.Range('E13').Value = '  +2.17%  '
$c = $ws.Range('D14')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '35.01'
$c.Style = $origStyle
$ws.Range('E14').Value = '  +0.26%  '
$c = $ws.Range('D15')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.655.69'
$c.Style = $origStyle
$ws.Range('E15').Value = '  +0.82%  '
$ws.Range('E16').Value = '  +2.87%  '
$c = $ws.Range('D17')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '64.297.40'
$c.Style = $origStyle
$ws.Range('E17').Value = '  +1.07%  '
$c = $ws.Range('D18')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.149.74'
$c.Style = $origStyle
$ws.Range('E18').Value = '  +1.08%  '
$c = $ws.Range('D19')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '6.85'
$c.Style = $origStyle
$ws.Range('E19').Value = '  +1.34%  '
$c = $ws.Range('D20')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '480.12'
$c.Style = $origStyle
$ws.Range('E20').Value = '  +1.32%  '
$c = $ws.Range('D21')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '14.51'
$c.Style = $origStyle
$ws.Range('E21').Value = '  +0.02%  '
$c = $ws.Range('D22')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.707'
$c.Style = $origStyle
$ws.Range('E22').Value = '  +0.29%  '
$c = $ws.Range('D23')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '7.65'
$c.Style = $origStyle
$ws.Range('E23').Value = '  -0.06%  '
$c = $ws.Range('D24')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '86.94'
$c.Style = $origStyle
$ws.Range('E24').Value = '  +3.72%  '
$c = $ws.Range('D25')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '13.35'
$c.Style = $origStyle
$ws.Range('E25').Value = '  -0.63%  '
$ws.Range('E26').Value = '  -0.01%  '
$c = $ws.Range('D27')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.74'
$c.Style = $origStyle
$ws.Range('E27').Value = '  -0.80%  '
$c = $ws.Range('D28')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '8.29'
$c.Style = $origStyle
$ws.Range('E28').Value = '  -0.50%  '
$c = $ws.Range('D29')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '7.15'
$c.Style = $origStyle
$ws.Range('E29').Value = '  +5.62%  '
$c = $ws.Range('D30')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.05'
$c.Style = $origStyle
$ws.Range('E30').Value = '  -1.59%  '
$c = $ws.Range('D31')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.112'
$c.Style = $origStyle
$ws.Range('E31').Value = '  -0.14%  '
$ws.Range('E32').Value = '  +0.09%  '
$c = $ws.Range('D33')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '26.67'
$c.Style = $origStyle
$ws.Range('E33').Value = '  +2.36%  '
$c = $ws.Range('D34')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.66'
$c.Style = $origStyle
$ws.Range('E34').Value = '  +0.13%  '
$ws.Range('E35').Value = '  -0.91%  '
$c = $ws.Range('D36')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '5.98'
$c.Style = $origStyle
$ws.Range('E36').Value = '  +1.52%  '
$c = $ws.Range('D37')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.0₃0754'
$c.Style = $origStyle
$ws.Range('E37').Value = '  +1.50%  '
$c = $ws.Range('D38')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '52.53'
$c.Style = $origStyle
$ws.Range('E38').Value = '  -0.28%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Range('D39')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '442.39'
$c.Style = $origStyle
$ws.Range('E39').Value = '  -2.20%  '
$ws.Range('B40').Value = 'dogwifhat'
$ws.Range('C40').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range('D40')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.94'
$c.Style = $origStyle
$ws.Range('E40').Value = '  +1.40%  '
$c = $ws.Range('D41')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.0391'
$c.Style = $origStyle
$ws.Range('E41').Value = '  +0.61%  '
$ws.Range('E42').Value = '  +1.59%  '
$c = $ws.Range('D43')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '8.20'
$c.Style = $origStyle
$ws.Range('E43').Value = '  -0.72%  '
$c = $ws.Range('D44')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.866.17'
$c.Style = $origStyle
$ws.Range('E44').Value = '  +1.23%  '
$c = $ws.Range('D45')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.259'
$c.Style = $origStyle
$ws.Range('E45').Value = '  -1.18%  '
$c = $ws.Range('D46')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.22'
$c.Style = $origStyle
$ws.Range('E46').Value = '  -0.99%  '
$c = $ws.Range('D47')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.44'
$c.Style = $origStyle
$ws.Range('E47').Value = '  +1.48%  '
$ws.Range('E48').Value = '  -0.03%  '
$c = $ws.Range('D49')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '25.89'
$c.Style = $origStyle
$ws.Range('E49').Value = '  +0.37%  '
$c = $ws.Range('D50')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.113'
$c.Style = $origStyle
$ws.Range('E50').Value = '  +0.67%  '
$c = $ws.Range('D51')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '121.59'
$c.Style = $origStyle
$ws.Range('E51').Value = '  +3.49%  '
